$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 25 (Sandringham Line), shifting the
# existing rows 25-27 down to 26-28.
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the "Richmond" entry.
$ws.Range("A25").Value = "Richmond"
$ws.Range("B25").Value = "MCG Zone 5 - The Great Southern Stand, Level 4, Bay Q18, Brunton Avenue"
$ws.Range("C25").Value = "27/12/20 12:30pm - 3:30pm"
$ws.Range("D25").Value = "Case did not attend during infectious period but may have acquired their illness here"

# Update the site name for the (now) last Springvale row.
$ws.Range("B28").Value = "Springvale Central,  268 Springvale Road"
